$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205262/1743535585915/aNjCWRFZTU-camera-video-segment-driver-1743535590915.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFKIZRCMV%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T200023Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJIMEYCIQD2bkzMtij0b8VixxnoMIF5OwgCs6jOqGFDz4u4EiRdawIhAPhlzeNtG%2BqxhZNp3IugpaQq4FeCoqzRpBQvzYCj88mRKuYDCNv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgzwzEst%2BqM7UQx70oIqugNCcQbaf8RWvnO0b%2FzG8lQJc3r5VSqE4DDvGryCSK%2F5hODqCyYFxJcWeZcB0TUOfQlIrcxR%2FXAZs3lFm5vRKmG4UVOdiBsSEQyrNJQRg2qRJVwQd3X21w7pAdNHZCWcuGHEpB9VdYhDHc4qKQ%2Bmin4IDDl01Z1%2ByUJRkC1dznYehPJiuQNdTxfBZnEJohOCjsVUKPIOSw6FMoa6aWjeZQnkYYrQEUorQRXshPt3T4Y1EOmUeN2RfMZwi1NQ4rl1JIcGYIK%2BTCH86ajwaUdDMk6%2BCZLbbyGXs%2FloRJCUz6ixDo1sbuXYEteL1ZJiZWpNvdgrWl2nW%2FUh2DpZL2i6JmgL%2BbD9i%2FpGs8CVl4NaaGNOarLDi6KRPMR8GS%2BdEYxJ8tdgs6v%2BiKImpVxLLc%2F7aEdeLSOMGbPI1LfelnAmU%2F95M9nyU9Fkd8jlM1tEbNnYNcaKD%2BMEF1hQvjiclz%2Fw%2FT5zM638LIk4%2BYqUhYci69NdEo0Q6evfGwXeYc6kMyqDyTSk6GH3DudpyFQz3QJS%2FVB1iq%2Bc2T9RWsZnUco%2FUhDCgxyOvrYid3Yq3g5BT%2FdOUdUPYt4nM17j5CR3MK7utb8GOqQBLJxOY6dOXgl4WKld%2FlHETJ6uPTlQAwa6jgCrJ%2BaI4JFNYbc5mfHtDcEeVYP9eTdtpDE1Y%2F2GqrMVtCIMnbNHXVOeHv3veJX326H9CLNtlErw2M%2BIoZsa5DedwkzduzOFcUcbgGiYMgnvSTYabwxpNKrSClOGa0MwPnys4BXS4krCHETKmUxEx8n22cTFLKMdMTsaw%2BQ6SfgzxJqGqviNWYmBsww%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2004%3A00%3A23%20GMT&X-Amz-Signature=9f455dbbe2c340397192b624f756c6d23a4e0d496fc8398b3ed8b3bbde4943e6'

$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535585915/3RxOLQEM5t-camera-video-segment-1743535590915.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFKIZRCMV%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T200023Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJIMEYCIQD2bkzMtij0b8VixxnoMIF5OwgCs6jOqGFDz4u4EiRdawIhAPhlzeNtG%2BqxhZNp3IugpaQq4FeCoqzRpBQvzYCj88mRKuYDCNv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgzwzEst%2BqM7UQx70oIqugNCcQbaf8RWvnO0b%2FzG8lQJc3r5VSqE4DDvGryCSK%2F5hODqCyYFxJcWeZcB0TUOfQlIrcxR%2FXAZs3lFm5vRKmG4UVOdiBsSEQyrNJQRg2qRJVwQd3X21w7pAdNHZCWcuGHEpB9VdYhDHc4qKQ%2Bmin4IDDl01Z1%2ByUJRkC1dznYehPJiuQNdTxfBZnEJohOCjsVUKPIOSw6FMoa6aWjeZQnkYYrQEUorQRXshPt3T4Y1EOmUeN2RfMZwi1NQ4rl1JIcGYIK%2BTCH86ajwaUdDMk6%2BCZLbbyGXs%2FloRJCUz6ixDo1sbuXYEteL1ZJiZWpNvdgrWl2nW%2FUh2DpZL2i6JmgL%2BbD9i%2FpGs8CVl4NaaGNOarLDi6KRPMR8GS%2BdEYxJ8tdgs6v%2BiKImpVxLLc%2F7aEdeLSOMGbPI1LfelnAmU%2F95M9nyU9Fkd8jlM1tEbNnYNcaKD%2BMEF1hQvjiclz%2Fw%2FT5zM638LIk4%2BYqUhYci69NdEo0Q6evfGwXeYc6kMyqDyTSk6GH3DudpyFQz3QJS%2FVB1iq%2Bc2T9RWsZnUco%2FUhDCgxyOvrYid3Yq3g5BT%2FdOUdUPYt4nM17j5CR3MK7utb8GOqQBLJxOY6dOXgl4WKld%2FlHETJ6uPTlQAwa6jgCrJ%2BaI4JFNYbc5mfHtDcEeVYP9eTdtpDE1Y%2F2GqrMVtCIMnbNHXVOeHv3veJX326H9CLNtlErw2M%2BIoZsa5DedwkzduzOFcUcbgGiYMgnvSTYabwxpNKrSClOGa0MwPnys4BXS4krCHETKmUxEx8n22cTFLKMdMTsaw%2BQ6SfgzxJqGqviNWYmBsww%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2004%3A00%3A23%20GMT&X-Amz-Signature=996f6ae68d0bf43188f82bdac35c0ffe4101fa039cb64a50ee414c666f717212'

$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535541388/Uf5dS3fS3z-camera-video-segment-driver-1743535543888.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFKIZRCMV%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T200023Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJIMEYCIQD2bkzMtij0b8VixxnoMIF5OwgCs6jOqGFDz4u4EiRdawIhAPhlzeNtG%2BqxhZNp3IugpaQq4FeCoqzRpBQvzYCj88mRKuYDCNv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgzwzEst%2BqM7UQx70oIqugNCcQbaf8RWvnO0b%2FzG8lQJc3r5VSqE4DDvGryCSK%2F5hODqCyYFxJcWeZcB0TUOfQlIrcxR%2FXAZs3lFm5vRKmG4UVOdiBsSEQyrNJQRg2qRJVwQd3X21w7pAdNHZCWcuGHEpB9VdYhDHc4qKQ%2Bmin4IDDl01Z1%2ByUJRkC1dznYehPJiuQNdTxfBZnEJohOCjsVUKPIOSw6FMoa6aWjeZQnkYYrQEUorQRXshPt3T4Y1EOmUeN2RfMZwi1NQ4rl1JIcGYIK%2BTCH86ajwaUdDMk6%2BCZLbbyGXs%2FloRJCUz6ixDo1sbuXYEteL1ZJiZWpNvdgrWl2nW%2FUh2DpZL2i6JmgL%2BbD9i%2FpGs8CVl4NaaGNOarLDi6KRPMR8GS%2BdEYxJ8tdgs6v%2BiKImpVxLLc%2F7aEdeLSOMGbPI1LfelnAmU%2F95M9nyU9Fkd8jlM1tEbNnYNcaKD%2BMEF1hQvjiclz%2Fw%2FT5zM638LIk4%2BYqUhYci69NdEo0Q6evfGwXeYc6kMyqDyTSk6GH3DudpyFQz3QJS%2FVB1iq%2Bc2T9RWsZnUco%2FUhDCgxyOvrYid3Yq3g5BT%2FdOUdUPYt4nM17j5CR3MK7utb8GOqQBLJxOY6dOXgl4WKld%2FlHETJ6uPTlQAwa6jgCrJ%2BaI4JFNYbc5mfHtDcEeVYP9eTdtpDE1Y%2F2GqrMVtCIMnbNHXVOeHv3veJX326H9CLNtlErw2M%2BIoZsa5DedwkzduzOFcUcbgGiYMgnvSTYabwxpNKrSClOGa0MwPnys4BXS4krCHETKmUxEx8n22cTFLKMdMTsaw%2BQ6SfgzxJqGqviNWYmBsww%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2004%3A00%3A23%20GMT&X-Amz-Signature=920bdafd7f188e7b768649ae26672cfb2aac529de3e458f731983d2ea6c2a408'

$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532925726/R4KzSHnOTc-camera-video-segment-driver-1743532928226.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFKIZRCMV%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T200023Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJIMEYCIQD2bkzMtij0b8VixxnoMIF5OwgCs6jOqGFDz4u4EiRdawIhAPhlzeNtG%2BqxhZNp3IugpaQq4FeCoqzRpBQvzYCj88mRKuYDCNv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgzwzEst%2BqM7UQx70oIqugNCcQbaf8RWvnO0b%2FzG8lQJc3r5VSqE4DDvGryCSK%2F5hODqCyYFxJcWeZcB0TUOfQlIrcxR%2FXAZs3lFm5vRKmG4UVOdiBsSEQyrNJQRg2qRJVwQd3X21w7pAdNHZCWcuGHEpB9VdYhDHc4qKQ%2Bmin4IDDl01Z1%2ByUJRkC1dznYehPJiuQNdTxfBZnEJohOCjsVUKPIOSw6FMoa6aWjeZQnkYYrQEUorQRXshPt3T4Y1EOmUeN2RfMZwi1NQ4rl1JIcGYIK%2BTCH86ajwaUdDMk6%2BCZLbbyGXs%2FloRJCUz6ixDo1sbuXYEteL1ZJiZWpNvdgrWl2nW%2FUh2DpZL2i6JmgL%2BbD9i%2FpGs8CVl4NaaGNOarLDi6KRPMR8GS%2BdEYxJ8tdgs6v%2BiKImpVxLLc%2F7aEdeLSOMGbPI1LfelnAmU%2F95M9nyU9Fkd8jlM1tEbNnYNcaKD%2BMEF1hQvjiclz%2Fw%2FT5zM638LIk4%2BYqUhYci69NdEo0Q6evfGwXeYc6kMyqDyTSk6GH3DudpyFQz3QJS%2FVB1iq%2Bc2T9RWsZnUco%2FUhDCgxyOvrYid3Yq3g5BT%2FdOUdUPYt4nM17j5CR3MK7utb8GOqQBLJxOY6dOXgl4WKld%2FlHETJ6uPTlQAwa6jgCrJ%2BaI4JFNYbc5mfHtDcEeVYP9eTdtpDE1Y%2F2GqrMVtCIMnbNHXVOeHv3veJX326H9CLNtlErw2M%2BIoZsa5DedwkzduzOFcUcbgGiYMgnvSTYabwxpNKrSClOGa0MwPnys4BXS4krCHETKmUxEx8n22cTFLKMdMTsaw%2BQ6SfgzxJqGqviNWYmBsww%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2004%3A00%3A23%20GMT&X-Amz-Signature=562a42ae7554262f688e7d5b8958b651f416080ac6ae2787139cac00bb6cfaba'

$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532886131/8jJIZWs74D-camera-video-segment-driver-1743532888631.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSFKIZRCMV%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T200023Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHIaCXVzLXdlc3QtMiJIMEYCIQD2bkzMtij0b8VixxnoMIF5OwgCs6jOqGFDz4u4EiRdawIhAPhlzeNtG%2BqxhZNp3IugpaQq4FeCoqzRpBQvzYCj88mRKuYDCNv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEQBBoMNzgxMjA0OTQyMjQ0IgzwzEst%2BqM7UQx70oIqugNCcQbaf8RWvnO0b%2FzG8lQJc3r5VSqE4DDvGryCSK%2F5hODqCyYFxJcWeZcB0TUOfQlIrcxR%2FXAZs3lFm5vRKmG4UVOdiBsSEQyrNJQRg2qRJVwQd3X21w7pAdNHZCWcuGHEpB9VdYhDHc4qKQ%2Bmin4IDDl01Z1%2ByUJRkC1dznYehPJiuQNdTxfBZnEJohOCjsVUKPIOSw6FMoa6aWjeZQnkYYrQEUorQRXshPt3T4Y1EOmUeN2RfMZwi1NQ4rl1JIcGYIK%2BTCH86ajwaUdDMk6%2BCZLbbyGXs%2FloRJCUz6ixDo1sbuXYEteL1ZJiZWpNvdgrWl2nW%2FUh2DpZL2i6JmgL%2BbD9i%2FpGs8CVl4NaaGNOarLDi6KRPMR8GS%2BdEYxJ8tdgs6v%2BiKImpVxLLc%2F7aEdeLSOMGbPI1LfelnAmU%2F95M9nyU9Fkd8jlM1tEbNnYNcaKD%2BMEF1hQvjiclz%2Fw%2FT5zM638LIk4%2BYqUhYci69NdEo0Q6evfGwXeYc6kMyqDyTSk6GH3DudpyFQz3QJS%2FVB1iq%2Bc2T9RWsZnUco%2FUhDCgxyOvrYid3Yq3g5BT%2FdOUdUPYt4nM17j5CR3MK7utb8GOqQBLJxOY6dOXgl4WKld%2FlHETJ6uPTlQAwa6jgCrJ%2BaI4JFNYbc5mfHtDcEeVYP9eTdtpDE1Y%2F2GqrMVtCIMnbNHXVOeHv3veJX326H9CLNtlErw2M%2BIoZsa5DedwkzduzOFcUcbgGiYMgnvSTYabwxpNKrSClOGa0MwPnys4BXS4krCHETKmUxEx8n22cTFLKMdMTsaw%2BQ6SfgzxJqGqviNWYmBsww%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2004%3A00%3A23%20GMT&X-Amz-Signature=5e80ec958b309d72b888b6d35ba0f80ef8319134679f67bce20eb6ade414e594'
